$wb = $excel.ActiveWorkbook

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutWs = $wb.Worksheets.Item("About")
$aboutWs.Range("A2").Value = "Version: $newVersion"
$aboutWs.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Jinyang Coal Mine, China, M2944, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 7; $r++) {
    $dataWs.Cells.Item($r, 19).Value = $newVersion
}
